$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semilla 8")

$ws.Range("C12").Value = "3045981670"
$ws.Range("B9").Value = "459399130"
$ws.Range("B11").Value = "255188531"
$ws.Range("B12").Value = "194936717"
$ws.Range("B13").Value = "432694001"
$ws.Range("B10").Value = "836898669"
$ws.Range("C13").Value = "3045981684"
$ws.Range("B14").Value = 920626579

$ws.Activate() | Out-Null
$ws.Range("B10").Select() | Out-Null
